$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptocurrency price/volume data and refreshed hour (G column: 9 -> 10)
$data = @(
    @{ Row = 2; D = "278.14"; E = "0.91%"; G = "10" },
    @{ Row = 3; D = "27.27"; E = "1.95%"; G = "10" },
    @{ Row = 4; D = "4.869"; E = "-0.11%"; G = "10" },
    @{ Row = 5; D = "0.06436"; E = "1.74%"; G = "10" },
    @{ Row = 6; D = "7.000"; E = "1.22%"; G = "10" },
    @{ Row = 7; D = "1.188"; E = "-7.20%"; G = "10" },
    @{ Row = 8; D = "0.8866"; E = "1.28%"; G = "10" },
    @{ Row = 9; D = "0.1571"; E = "2.76%"; G = "10" },
    @{ Row = 10; D = "0.05112"; E = "2.32%"; G = "10" },
    @{ Row = 11; D = "0.07505"; E = "0.36%"; G = "10" },
    @{ Row = 12; D = "0.02884"; E = "-3.86%"; G = "10" },
    @{ Row = 13; D = "0.08971"; E = "-0.89%"; G = "10" },
    @{ Row = 14; D = "0.001571"; E = "-0.39%"; G = "10" },
    @{ Row = 15; D = "0.0006367"; E = "0.57%"; G = "10" },
    @{ Row = 16; D = "0.006141"; E = "4.50%"; G = "10" },
    @{ Row = 17; D = "3.479"; G = "10" },
    @{ Row = 18; D = "3.307"; E = "-0.44%"; G = "10" },
    @{ Row = 19; D = "2.273"; E = "0.06%"; G = "10" },
    @{ Row = 20; E = "1.10%"; G = "10" },
    @{ Row = 21; E = "1.10%"; G = "10" },
    @{ Row = 22; D = "3.931"; E = "0.73%"; G = "10" },
    @{ Row = 23; D = "0.04413"; E = "1.33%"; G = "10" },
    @{ Row = 24; G = "10" },
    @{ Row = 25; D = "0.001175"; E = "0.38%"; G = "10" },
    @{ Row = 26; D = "0.003874"; E = "-8.04%"; G = "10" },
    @{ Row = 27; G = "10" },
    @{ Row = 28; E = "-1.70%"; G = "10" },
    @{ Row = 29; E = "-1.74%"; G = "10" },
    @{ Row = 30; G = "10" },
    @{ Row = 31; G = "10" },
    @{ Row = 32; G = "10" },
    @{ Row = 33; G = "10" },
    @{ Row = 34; G = "10" },
    @{ Row = 35; G = "10" },
    @{ Row = 36; G = "10" },
    @{ Row = 37; G = "10" },
    @{ Row = 38; G = "10" },
    @{ Row = 39; G = "10" },
    @{ Row = 40; D = "0.04138"; E = "0.82%"; G = "10" },
    @{ Row = 41; D = "0.006737"; E = "-3.42%"; G = "10" },
    @{ Row = 42; E = "0.21%"; G = "10" },
    @{ Row = 43; D = "0.001870"; E = "-12.78%"; G = "10" },
    @{ Row = 44; D = "0.01122"; E = "3.89%"; G = "10" },
    @{ Row = 45; D = "0.00005304"; E = "0.37%"; G = "10" },
    @{ Row = 46; G = "10" },
    @{ Row = 47; D = "0.01852"; E = "-11.82%"; G = "10" },
    @{ Row = 48; G = "10" },
    @{ Row = 49; G = "10" },
    @{ Row = 50; G = "10" },
    @{ Row = 51; G = "10" }
)

foreach ($item in $data) {
    $r = $item.Row
    if ($item.ContainsKey("D")) {
        $ws.Cells.Item($r, 4).Value = "'" + $item.D
    }
    if ($item.ContainsKey("E")) {
        $ws.Cells.Item($r, 5).Value = "'" + $item.E
    }
    if ($item.ContainsKey("G")) {
        $ws.Cells.Item($r, 7).Value = "'" + $item.G
    }
}
